$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.4394200868768915
$ws.Cells.Item(2, 3).Value = 0.7311110424551147
$ws.Cells.Item(2, 4).Value = 1.682578117590849
$ws.Cells.Item(2, 5).Value = 1.297142288876147
$ws.Cells.Item(2, 6).Value = 1.235239686460398
$ws.Cells.Item(2, 7).Value = 42
$ws.Cells.Item(3, 2).Value = 0.2153962797860198
$ws.Cells.Item(3, 3).Value = 1.302986795684886
$ws.Cells.Item(3, 4).Value = 4.677571906077256
$ws.Cells.Item(3, 5).Value = 2.162769499063008
$ws.Cells.Item(3, 6).Value = 2.178750962696262
$ws.Cells.Item(3, 7).Value = 41
$ws.Cells.Item(4, 2).Value = 0.504178189794622
$ws.Cells.Item(4, 3).Value = 1.31039634774812
$ws.Cells.Item(4, 4).Value = 4.184774737970552
$ws.Cells.Item(4, 5).Value = 2.04567219709575
$ws.Cells.Item(4, 6).Value = 2.007825482993971
$ws.Cells.Item(4, 7).Value = 40
$ws.Cells.Item(5, 2).Value = 0.322232619466939
$ws.Cells.Item(5, 3).Value = 1.251316763704151
$ws.Cells.Item(5, 4).Value = 3.808272273858065
$ws.Cells.Item(5, 5).Value = 1.951479508951622
$ws.Cells.Item(5, 6).Value = 1.949852208296635
$ws.Cells.Item(5, 7).Value = 39
$ws.Cells.Item(6, 2).Value = 0.5079447904278251
$ws.Cells.Item(6, 3).Value = 1.466392046323584
$ws.Cells.Item(6, 4).Value = 4.538018304700308
$ws.Cells.Item(6, 5).Value = 2.130262496665683
$ws.Cells.Item(6, 6).Value = 2.096589218513666
$ws.Cells.Item(6, 7).Value = 38
$ws.Cells.Item(7, 2).Value = 0.3328774944335463
$ws.Cells.Item(7, 3).Value = 1.415444145779199
$ws.Cells.Item(7, 4).Value = 4.344522175199605
$ws.Cells.Item(7, 5).Value = 2.084351739798157
$ws.Cells.Item(7, 6).Value = 2.085981288595052
$ws.Cells.Item(7, 7).Value = 37
$ws.Cells.Item(8, 2).Value = 0.5270938661561243
$ws.Cells.Item(8, 3).Value = 1.466508675359905
$ws.Cells.Item(8, 4).Value = 4.698901303729951
$ws.Cells.Item(8, 5).Value = 2.167694928658078
$ws.Cells.Item(8, 6).Value = 2.132460959010634
$ws.Cells.Item(8, 7).Value = 36
$ws.Cells.Item(9, 2).Value = 0.2983617769427525
$ws.Cells.Item(9, 3).Value = 1.395751893831891
$ws.Cells.Item(9, 4).Value = 4.30447533330254
$ws.Cells.Item(9, 5).Value = 2.074722953384991
$ws.Cells.Item(9, 6).Value = 2.08313215401905
$ws.Cells.Item(9, 7).Value = 35
$ws.Cells.Item(10, 2).Value = 0.4106355216869119
$ws.Cells.Item(10, 3).Value = 1.33181755496485
$ws.Cells.Item(10, 4).Value = 4.092714757033598
$ws.Cells.Item(10, 5).Value = 2.023045910757736
$ws.Cells.Item(10, 6).Value = 2.010722542093411
$ws.Cells.Item(10, 7).Value = 34
$ws.Cells.Item(11, 2).Value = 0.3748234707463745
$ws.Cells.Item(11, 3).Value = 1.370250310854823
$ws.Cells.Item(11, 4).Value = 4.504998655405332
$ws.Cells.Item(11, 5).Value = 2.122498210931008
$ws.Cells.Item(11, 6).Value = 2.121531718910878
$ws.Cells.Item(11, 7).Value = 33
